$d = $word.ActiveDocument
$r = $d.Content
$r.Find.Execute("20+61=81", $true, $false, $false, $false, $false, $true, 1, $false, "37+20=57", 1) | Out-Null
$r.Find.Execute("69-31=38", $true, $false, $false, $false, $false, $true, 1, $false, "79-36=43", 1) | Out-Null
$r.Find.Execute("99-59=40", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=43", 1) | Out-Null
$r.Find.Execute("12+75=87", $true, $false, $false, $false, $false, $true, 1, $false, "13+84=97", 1) | Out-Null
$r.Find.Execute("27+41=68", $true, $false, $false, $false, $false, $true, 1, $false, "60+22=82", 1) | Out-Null
$r.Find.Execute("83+5=88", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=36", 1) | Out-Null
$r.Find.Execute("98-70=28", $true, $false, $false, $false, $false, $true, 1, $false, "10+20=30", 1) | Out-Null
$r.Find.Execute("50-10=40", $true, $false, $false, $false, $false, $true, 1, $false, "66+18=84", 1) | Out-Null
$r.Find.Execute("51+27=78", $true, $false, $false, $false, $false, $true, 1, $false, "1+33=34", 1) | Out-Null
$r.Find.Execute("36+3=39", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=1", 1) | Out-Null
$r.Find.Execute("67-2=65", $true, $false, $false, $false, $false, $true, 1, $false, "86-76=10", 1) | Out-Null
$r.Find.Execute("34+43=77", $true, $false, $false, $false, $false, $true, 1, $false, "84-4=80", 1) | Out-Null
$r.Find.Execute("21+77=98", $true, $false, $false, $false, $false, $true, 1, $false, "38-9=29", 1) | Out-Null
$r.Find.Execute("71-71=0", $true, $false, $false, $false, $false, $true, 1, $false, "82-60=22", 1) | Out-Null
$r.Find.Execute("26-1=25", $true, $false, $false, $false, $false, $true, 1, $false, "20-12=8", 1) | Out-Null
$r.Find.Execute("68+9=77", $true, $false, $false, $false, $false, $true, 1, $false, "61-9=52", 1) | Out-Null
$r.Find.Execute("26+0=26", $true, $false, $false, $false, $false, $true, 1, $false, "89-22=67", 1) | Out-Null
$r.Find.Execute("66+1=67", $true, $false, $false, $false, $false, $true, 1, $false, "55-23=32", 1) | Out-Null
$r.Find.Execute("48-24=24", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=55", 1) | Out-Null
$r.Find.Execute("50+1=51", $true, $false, $false, $false, $false, $true, 1, $false, "44+34=78", 1) | Out-Null
$r.Find.Execute("34+31=65", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=67", 1) | Out-Null
$r.Find.Execute("26+42=68", $true, $false, $false, $false, $false, $true, 1, $false, "52+44=96", 1) | Out-Null
$r.Find.Execute("60-56=4", $true, $false, $false, $false, $false, $true, 1, $false, "94-54=40", 1) | Out-Null
$r.Find.Execute("53+45=98", $true, $false, $false, $false, $false, $true, 1, $false, "69-14=55", 1) | Out-Null
$r.Find.Execute("81-2=79", $true, $false, $false, $false, $false, $true, 1, $false, "91-81=10", 1) | Out-Null
$r.Find.Execute("92-59=33", $true, $false, $false, $false, $false, $true, 1, $false, "42+51=93", 1) | Out-Null
$r.Find.Execute("29+31=60", $true, $false, $false, $false, $false, $true, 1, $false, "16+3=19", 1) | Out-Null
$r.Find.Execute("37+3=40", $true, $false, $false, $false, $false, $true, 1, $false, "17+72=89", 1) | Out-Null
$r.Find.Execute("18+52=70", $true, $false, $false, $false, $false, $true, 1, $false, "63+0=63", 1) | Out-Null
$r.Find.Execute("65-56=9", $true, $false, $false, $false, $false, $true, 1, $false, "40-3=37", 1) | Out-Null
$r.Find.Execute("15+59=74", $true, $false, $false, $false, $false, $true, 1, $false, "5+92=97", 1) | Out-Null
$r.Find.Execute("97-39=58", $true, $false, $false, $false, $false, $true, 1, $false, "14-2=12", 1) | Out-Null
$r.Find.Execute("86-42=44", $true, $false, $false, $false, $false, $true, 1, $false, "12+32=44", 1) | Out-Null
$r.Find.Execute("18+48=66", $true, $false, $false, $false, $false, $true, 1, $false, "25+7=32", 1) | Out-Null
$r.Find.Execute("60+4=64", $true, $false, $false, $false, $false, $true, 1, $false, "63-53=10", 1) | Out-Null
$r.Find.Execute("41+45=86", $true, $false, $false, $false, $false, $true, 1, $false, "93-76=17", 1) | Out-Null
$r.Find.Execute("26+7=33", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=32", 1) | Out-Null
$r.Find.Execute("23+43=66", $true, $false, $false, $false, $false, $true, 1, $false, "9+6=15", 1) | Out-Null
$r.Find.Execute("67-2=65", $true, $false, $false, $false, $false, $true, 1, $false, "52+43=95", 1) | Out-Null
$r.Find.Execute("30+64=94", $true, $false, $false, $false, $false, $true, 1, $false, "12+41=53", 1) | Out-Null
$r.Find.Execute("25-20=5", $true, $false, $false, $false, $false, $true, 1, $false, "45+20=65", 1) | Out-Null
$r.Find.Execute("36+15=51", $true, $false, $false, $false, $false, $true, 1, $false, "86+11=97", 1) | Out-Null
$r.Find.Execute("72+7=79", $true, $false, $false, $false, $false, $true, 1, $false, "4+31=35", 1) | Out-Null
$r.Find.Execute("91-79=12", $true, $false, $false, $false, $false, $true, 1, $false, "59-20=39", 1) | Out-Null
$r.Find.Execute("24+6=30", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=22", 1) | Out-Null
$r.Find.Execute("9+16=25", $true, $false, $false, $false, $false, $true, 1, $false, "63-19=44", 1) | Out-Null
$r.Find.Execute("55+6=61", $true, $false, $false, $false, $false, $true, 1, $false, "58+39=97", 1) | Out-Null
$r.Find.Execute("45+10=55", $true, $false, $false, $false, $false, $true, 1, $false, "65+18=83", 1) | Out-Null
$r.Find.Execute("21-14=7", $true, $false, $false, $false, $false, $true, 1, $false, "22+58=80", 1) | Out-Null
$r.Find.Execute("19+61=80", $true, $false, $false, $false, $false, $true, 1, $false, "73-6=67", 1) | Out-Null
$r.Find.Execute("82+13=95", $true, $false, $false, $false, $false, $true, 1, $false, "54+3=57", 1) | Out-Null
$r.Find.Execute("95-83=12", $true, $false, $false, $false, $false, $true, 1, $false, "55+39=94", 1) | Out-Null
$r.Find.Execute("26+19=45", $true, $false, $false, $false, $false, $true, 1, $false, "69-16=53", 1) | Out-Null
$r.Find.Execute("12+77=89", $true, $false, $false, $false, $false, $true, 1, $false, "13+43=56", 1) | Out-Null
$r.Find.Execute("29+17=46", $true, $false, $false, $false, $false, $true, 1, $false, "87-53=34", 1) | Out-Null
$r.Find.Execute("55-22=33", $true, $false, $false, $false, $false, $true, 1, $false, "70-37=33", 1) | Out-Null
$r.Find.Execute("19+13=32", $true, $false, $false, $false, $false, $true, 1, $false, "43+47=90", 1) | Out-Null
$r.Find.Execute("90-59=31", $true, $false, $false, $false, $false, $true, 1, $false, "50+28=78", 1) | Out-Null
$r.Find.Execute("11+35=46", $true, $false, $false, $false, $false, $true, 1, $false, "59+26=85", 1) | Out-Null
$r.Find.Execute("32+56=88", $true, $false, $false, $false, $false, $true, 1, $false, "42+37=79", 1) | Out-Null
$r.Find.Execute("20+42=62", $true, $false, $false, $false, $false, $true, 1, $false, "84-37=47", 1) | Out-Null
$r.Find.Execute("16+68=84", $true, $false, $false, $false, $false, $true, 1, $false, "49-11=38", 1) | Out-Null
$r.Find.Execute("29+41=70", $true, $false, $false, $false, $false, $true, 1, $false, "71+7=78", 1) | Out-Null
$r.Find.Execute("14+30=44", $true, $false, $false, $false, $false, $true, 1, $false, "7+88=95", 1) | Out-Null
$r.Find.Execute("32+61=93", $true, $false, $false, $false, $false, $true, 1, $false, "26+51=77", 1) | Out-Null
$r.Find.Execute("79+14=93", $true, $false, $false, $false, $false, $true, 1, $false, "61+30=91", 1) | Out-Null
$r.Find.Execute("11+16=27", $true, $false, $false, $false, $false, $true, 1, $false, "75-43=32", 1) | Out-Null
$r.Find.Execute("87-25=62", $true, $false, $false, $false, $false, $true, 1, $false, "18-8=10", 1) | Out-Null
$r.Find.Execute("95-4=91", $true, $false, $false, $false, $false, $true, 1, $false, "51+46=97", 1) | Out-Null
$r.Find.Execute("44-37=7", $true, $false, $false, $false, $false, $true, 1, $false, "53-26=27", 1) | Out-Null
$r.Find.Execute("46+12=58", $true, $false, $false, $false, $false, $true, 1, $false, "83-72=11", 1) | Out-Null
$r.Find.Execute("89-8=81", $true, $false, $false, $false, $false, $true, 1, $false, "96-80=16", 1) | Out-Null
$r.Find.Execute("87-50=37", $true, $false, $false, $false, $false, $true, 1, $false, "31+10=41", 1) | Out-Null
$r.Find.Execute("16+54=70", $true, $false, $false, $false, $false, $true, 1, $false, "0+9=9", 1) | Out-Null
$r.Find.Execute("82-77=5", $true, $false, $false, $false, $false, $true, 1, $false, "83-54=29", 1) | Out-Null
$r.Find.Execute("61+10=71", $true, $false, $false, $false, $false, $true, 1, $false, "2+30=32", 1) | Out-Null
$r.Find.Execute("72-47=25", $true, $false, $false, $false, $false, $true, 1, $false, "90-29=61", 1) | Out-Null
$r.Find.Execute("85-17=68", $true, $false, $false, $false, $false, $true, 1, $false, "16+1=17", 1) | Out-Null
$r.Find.Execute("51+32=83", $true, $false, $false, $false, $false, $true, 1, $false, "44-3=41", 1) | Out-Null
$r.Find.Execute("26-0=26", $true, $false, $false, $false, $false, $true, 1, $false, "90-49=41", 1) | Out-Null
$r.Find.Execute("32-22=10", $true, $false, $false, $false, $false, $true, 1, $false, "34+24=58", 1) | Out-Null
$r.Find.Execute("93-56=37", $true, $false, $false, $false, $false, $true, 1, $false, "88-37=51", 1) | Out-Null
$r.Find.Execute("13+71=84", $true, $false, $false, $false, $false, $true, 1, $false, "13-8=5", 1) | Out-Null
$r.Find.Execute("68-8=60", $true, $false, $false, $false, $false, $true, 1, $false, "85+3=88", 1) | Out-Null
$r.Find.Execute("54-37=17", $true, $false, $false, $false, $false, $true, 1, $false, "21+30=51", 1) | Out-Null
$r.Find.Execute("49-0=49", $true, $false, $false, $false, $false, $true, 1, $false, "28+47=75", 1) | Out-Null
$r.Find.Execute("83-82=1", $true, $false, $false, $false, $false, $true, 1, $false, "93-51=42", 1) | Out-Null
$r.Find.Execute("7+16=23", $true, $false, $false, $false, $false, $true, 1, $false, "26+67=93", 1) | Out-Null
$r.Find.Execute("68-56=12", $true, $false, $false, $false, $false, $true, 1, $false, "23+63=86", 1) | Out-Null
$r.Find.Execute("71-47=24", $true, $false, $false, $false, $false, $true, 1, $false, "28-20=8", 1) | Out-Null
$r.Find.Execute("76-2=74", $true, $false, $false, $false, $false, $true, 1, $false, "17+1=18", 1) | Out-Null
$r.Find.Execute("73+13=86", $true, $false, $false, $false, $false, $true, 1, $false, "27+37=64", 1) | Out-Null
$r.Find.Execute("30+53=83", $true, $false, $false, $false, $false, $true, 1, $false, "65+22=87", 1) | Out-Null
$r.Find.Execute("30+6=36", $true, $false, $false, $false, $false, $true, 1, $false, "76-40=36", 1) | Out-Null
$r.Find.Execute("90+3=93", $true, $false, $false, $false, $false, $true, 1, $false, "51+28=79", 1) | Out-Null
$r.Find.Execute("5+32=37", $true, $false, $false, $false, $false, $true, 1, $false, "62+20=82", 1) | Out-Null
$r.Find.Execute("80-30=50", $true, $false, $false, $false, $false, $true, 1, $false, "8+24=32", 1) | Out-Null
$r.Find.Execute("66-23=43", $true, $false, $false, $false, $false, $true, 1, $false, "9+65=74", 1) | Out-Null
$r.Find.Execute("96-3=93", $true, $false, $false, $false, $false, $true, 1, $false, "79+9=88", 1) | Out-Null
$r.Find.Execute("22+0=22", $true, $false, $false, $false, $false, $true, 1, $false, "53+28=81", 1) | Out-Null
